# Update "想去人数" (column F) and "最低票价" (column G) values
# on the "展览" and "全部类型" worksheets, matching the latest
# bilibili-scraped data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 245

    $ws.Range("F6").Value = 255

    $ws.Range("F7").Value = 6056
    $ws.Range("G7").Value = 58.5

    $ws.Range("F8").Value = 43

    $ws.Range("F11").Value = 57

    $ws.Range("F14").Value = 185

    $ws.Range("F15").Value = 406
}
